# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Sat Jul 22 16:56:30 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.901.12"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.889.73"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("D4").Value = "'0.9999"

$ws.Range("D5").Value = "'0.7685"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").Value = "'242.74"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.3132"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "'25.68"
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("D10").Value = "'0.07139"
$ws.Range("E10").Value = "  -2.50%  "

$ws.Range("D11").Value = "'0.08525"
$ws.Range("E11").Value = "  +4.85%  "

$ws.Range("D12").Value = "'0.7636"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.898.69"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.367"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").Value = "'93.80"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").Value = "'6.168"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "29.783.94"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").Value = "'244.22"
$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("D20").Value = "'0.000007807"
$ws.Range("E20").Value = "  -0.58%  "

$ws.Range("D21").Value = "'0.9993"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "'7.999"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'0.1628"
$ws.Range("E24").Value = "  +2.92%  "

$ws.Range("D25").Value = "'9.415"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").Value = "'162.99"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").Value = "'18.78"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("D28").Value = "'2.038"
$ws.Range("E28").Value = "  +0.38%  "

$ws.Range("D29").Value = "'1.505"
$ws.Range("E29").Value = "  +4.00%  "

$ws.Range("D30").Value = "'1.540"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").Value = "'4.497"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").Value = "'4.117"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").Value = "'0.05454"
$ws.Range("E33").Value = "  -2.31%  "

$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").Value = "'0.7457"
$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("D37").Value = "'2.698"
$ws.Range("E37").Value = "  +2.50%  "

$ws.Range("E38").Value = "  +0.83%  "

$ws.Range("D39").Value = "'2.782"
$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("D40").Value = "'0.4475"
$ws.Range("E40").Value = "  +0.88%  "

$ws.Range("D41").Value = "1.102.57"
$ws.Range("E41").Value = "  -3.15%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'73.26"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.088"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("D44").Value = "'0.8535"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "'102.85"
$ws.Range("E46").Value = "  +0.90%  "

$ws.Range("D47").Value = "'1.868"
$ws.Range("E47").Value = "  -1.38%  "

$ws.Range("D48").Value = "'7.664"
$ws.Range("E48").Value = "  +2.45%  "

$ws.Range("D49").Value = "'3.082"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06088"
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "1.997.74"
$ws.Range("E51").Value = "  -1.89%  "
